$d = $word.ActiveDocument

# Replace the "Launcher - rocket launcher" unit entry with the new
# "Hydro- Launches an ice rocket" text (Hydro is now the launcher unit,
# firing "ice rockets").
$d.Content.Find.Execute("Launcher - rocket launcher", $true, $true, $false, $false, $false, $true, 1, $false, "Hydro- Launches an ice rocket", 2)
